# Apply the EMMOntoPy "onto_update_only_classes.xlsx" update:
#  - ImportedOntologies!A3: bump the imported EMMO version from 1.0.0-beta
#    (emmo-inferred-chemistry2.ttl) to 1.0.0-beta4 (emmo-inferred.ttl) and
#    turn the cell into a clickable hyperlink to that URL.
#  - Concepts!A4: rename the "Pattern" example concept to "SpecialPattern".
#  - Restore the view state left by the author after the edit: Concepts is
#    the active/selected sheet, with D13 selected; ImportedOntologies has
#    A12 selected; Metadata keeps its previous B20 selection.

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsImported = $wb.Worksheets.Item("ImportedOntologies")
$wsConcepts = $wb.Worksheets.Item("Concepts")

# --- ImportedOntologies: update the ontology URL and hyperlink it ---
$newUrl = "https://raw.githubusercontent.com/emmo-repo/emmo-repo.github.io/master/versions/1.0.0-beta4/emmo-inferred.ttl"
$cellA3 = $wsImported.Range("A3")
$cellA3.Value = $newUrl
$wsImported.Hyperlinks.Add($cellA3, $newUrl, "", "", $newUrl)

# --- Concepts: rename the Pattern concept ---
$wsConcepts.Range("A4").Value = "SpecialPattern"

# --- View state: selections on each sheet and the active tab ---
$wsMetadata.Range("B20").Select()

$wsImported.Activate()
$wsImported.Range("A12").Select()

$wsConcepts.Activate()
$wsConcepts.Range("D13").Select()
